$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.355.47"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.432.82"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'574.78"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'145.45"
$ws.Range("E6").Value = "  +7.08%  "
$ws.Range("D7").Value = "3.433.50"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").Value = "'7.67"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "4.020.86"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'28.03"
$ws.Range("E14").Value = "  +8.01%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "3.431.46"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "61.470.45"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +7.94%  "
$ws.Range("D20").Value = "'14.18"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "'395.12"
$ws.Range("E22").Value = "  +6.41%  "
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").Value = "'73.00"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").Value = "'0.995"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").Value = "3.570.93"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "'0.177"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "'7.60"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").Value = "'1.46"
$ws.Range("E33").Value = "  -6.36%  "
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("D38").Value = "3.460.13"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "'167.37"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").Value = "'0.0787"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").Value = "'26.86"
$ws.Range("E43").Value = "  +7.27%  "
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").Value = "'42.08"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "2.579.73"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("E51").Value = "  +2.61%  "
